# Refresh the crypto price/volume snapshot (scheduled GitHub Actions pull).
#
# The Price column (D) is stored as plain text in the sheet (e.g. "67.824.01",
# "0.613", "9.15"). Several of the new prices are syntactically valid numbers,
# so assigning them straight to Range.Value would make Excel auto-convert them
# to doubles (and mangle the display, e.g. "0.613" -> 0.61299999999999999).
# Prefixing those values with a leading apostrophe tells Excel to keep them as
# literal text, exactly like a user typing '0.613 into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.824.01"
$ws.Range("E2").Value = "  +9.11%  "
# Row 3
$ws.Range("D3").Value = "3.519.33"
$ws.Range("E3").Value = "  +11.02%  "
# Row 4
$ws.Range("E4").Value = "  +0.11%  "
# Row 5
$ws.Range("D5").Value = "'191.73"
$ws.Range("E5").Value = "  +12.44%  "
# Row 6
$ws.Range("D6").Value = "'557.25"
$ws.Range("E6").Value = "  +9.50%  "
# Row 7
$ws.Range("D7").Value = "3.524.84"
$ws.Range("E7").Value = "  +11.35%  "
# Row 8
$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  +4.57%  "
# Row 9
$ws.Range("E9").Value = "  +0.04%  "
# Row 10
$ws.Range("D10").Value = "'0.639"
$ws.Range("E10").Value = "  +8.47%  "
# Row 11
$ws.Range("D11").Value = "'57.24"
$ws.Range("E11").Value = "  +5.62%  "
# Row 12
$ws.Range("E12").Value = "  +17.21%  "
# Row 13
$ws.Range("D13").Value = "'0.0000277"
$ws.Range("E13").Value = "  +11.10%  "
# Row 14
$ws.Range("E14").Value = "  +8.04%  "
# Row 15
$ws.Range("D15").Value = "4.082.42"
$ws.Range("E15").Value = "  +10.94%  "
# Row 16
$ws.Range("D16").Value = "3.524.25"
$ws.Range("E16").Value = "  +11.23%  "
# Row 17
$ws.Range("D17").Value = "68.436.59"
$ws.Range("E17").Value = "  +10.43%  "
# Row 18
$ws.Range("E18").Value = "  +7.59%  "
# Row 19
$ws.Range("D19").Value = "'18.41"
$ws.Range("E19").Value = "  +9.35%  "
# Row 20
$ws.Range("D20").Value = "'11.93"
$ws.Range("E20").Value = "  +12.14%  "
# Row 21
$ws.Range("E21").Value = "  +8.47%  "
# Row 22
$ws.Range("D22").Value = "'409.47"
$ws.Range("E22").Value = "  +14.29%  "
# Row 23
$ws.Range("E23").Value = "  +9.09%  "
# Row 24
$ws.Range("D24").Value = "'84.80"
$ws.Range("E24").Value = "  +7.87%  "
# Row 25
$ws.Range("E25").Value = "  +9.67%  "
# Row 26
$ws.Range("E26").Value = "  +11.02%  "
# Row 27
$ws.Range("E27").Value = "  +12.84%  "
# Row 28
$ws.Range("E28").Value = "  -0.11%  "
# Row 29
$ws.Range("D29").Value = "'11.92"
$ws.Range("E29").Value = "  +8.66%  "
# Row 30
$ws.Range("E30").Value = "  +7.95%  "
# Row 31
$ws.Range("E31").Value = "  +10.49%  "
# Row 32
$ws.Range("D32").Value = "'682.28"
$ws.Range("E32").Value = "  +10.41%  "
# Row 33
$ws.Range("D33").Value = "'6.90"
$ws.Range("E33").Value = "  +8.35%  "
# Row 34
$ws.Range("E34").Value = "  +7.48%  "
# Row 35
$ws.Range("E35").Value = "  +9.97%  "
# Row 36
$ws.Range("D36").Value = "'60.58"
$ws.Range("E36").Value = "  +7.15%  "
# Row 37
$ws.Range("D37").Value = "0.0₃0848"
$ws.Range("E37").Value = "  +28.59%  "
# Row 38
$ws.Range("D38").Value = "'39.18"
$ws.Range("E38").Value = "  +8.50%  "
# Row 39
$ws.Range("D39").Value = "'0.404"
$ws.Range("E39").Value = "  +8.84%  "
# Row 40
$ws.Range("D40").Value = "'0.999"
# Row 41
$ws.Range("E41").Value = "  +26.32%  "
# Row 42
$ws.Range("D42").Value = "'0.133"
$ws.Range("E42").Value = "  +11.90%  "
# Row 43
$ws.Range("D43").Value = "'2.76"
$ws.Range("E43").Value = "  +16.54%  "
# Row 44
$ws.Range("E44").Value = "  +17.23%  "
# Row 45
$ws.Range("E45").Value = "  +0.26%  "
# Row 46
$ws.Range("D46").Value = "3.053.01"
$ws.Range("E46").Value = "  +9.68%  "
# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0424"
$ws.Range("E47").Value = "  +11.72%  "
# Row 48
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'9.15"
$ws.Range("E48").Value = "  +24.05%  "
# Row 49
$ws.Range("E49").Value = "  +7.05%  "
# Row 50
$ws.Range("D50").Value = "'3.21"
$ws.Range("E50").Value = "  +11.31%  "
# Row 51
$ws.Range("E51").Value = "  +8.51%  "
